$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C (Förändrad / last-changed date) for all existing data
#    rows (2..496) from 45172 to 45175.
$ws.Range("C2:C496").Value = 45175

# 2) Row 496 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(496).RowHeight = 15

# 3) Append new row 497.
$ws.Range("A497").Value = "A 40994-2023"
$ws.Range("B497").Value = 45173
$ws.Range("B497").NumberFormat = "YYYY-MM-DD"
$ws.Range("C497").Value = 45175
$ws.Range("C497").NumberFormat = "YYYY-MM-DD"
$ws.Range("D497").Value = "HALLANDS LÄN"
$ws.Range("E497").Value = "FALKENBERG"
$ws.Range("G497").Value = 3.4
$ws.Range("H497").Value = 0
$ws.Range("I497").Value = 0
$ws.Range("J497").Value = 0
$ws.Range("K497").Value = 0
$ws.Range("L497").Value = 0
$ws.Range("M497").Value = 0
$ws.Range("N497").Value = 0
$ws.Range("O497").Value = 0
$ws.Range("P497").Value = 0
$ws.Range("Q497").Value = 0
$ws.Range("R497").WrapText = $true
$ws.Rows.Item(497).RowHeight = 15

# 4) Append new row 498 (no explicit custom row height).
$ws.Range("A498").Value = "A 41242-2023"
$ws.Range("B498").Value = 45174
$ws.Range("B498").NumberFormat = "YYYY-MM-DD"
$ws.Range("C498").Value = 45175
$ws.Range("C498").NumberFormat = "YYYY-MM-DD"
$ws.Range("D498").Value = "HALLANDS LÄN"
$ws.Range("E498").Value = "FALKENBERG"
$ws.Range("G498").Value = 0.6
$ws.Range("H498").Value = 0
$ws.Range("I498").Value = 0
$ws.Range("J498").Value = 0
$ws.Range("K498").Value = 0
$ws.Range("L498").Value = 0
$ws.Range("M498").Value = 0
$ws.Range("N498").Value = 0
$ws.Range("O498").Value = 0
$ws.Range("P498").Value = 0
$ws.Range("Q498").Value = 0
$ws.Range("R498").WrapText = $true

Write-Output "done"
